$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.520731687545776
$ws.Range("B1").Value = 3.665806770324707
$ws.Range("C1").Value = 5.86833667755127
$ws.Range("D1").Value = 1.431361317634583
$ws.Range("E1").Value = 0.8366617560386658
